$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "42.706.20"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "2.306.95"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.21%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.994"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").Value = "2.656.42"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").Value = "2.301.81"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "42.867.28"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.95%  "

$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.43%  "

$ws.Range("E23").Value = "  -2.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +11.77%  "

$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.68"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("E33").Value = "  -2.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.81"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("D38").ClearFormats()

$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "108.34"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +10.77%  "

$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.93%  "

$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("D47").Value = "1.722.12"
$ws.Range("E47").Value = "  +4.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.66"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.17"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.48%  "
